# Atualização de bases das ligas, do dia: 01-06-2024 às 01:16
#
# The source data rows got re-ordered upstream; for five adjacent row
# pairs the "id" (column A) stayed put but every other field (columns
# B..AD: match id, Div, Date, HomeTeam, AwayTeam, scores, odds, ...)
# needs to swap between the two rows.
#
# Row pairs (1-based sheet rows) that swap their B:AD payload:
#   107 <-> 108
#   128 <-> 129
#   148 <-> 149
#   153 <-> 154
#   211 <-> 212

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Scratch cell well outside the used range (A1:AD261) used as a holding
# spot while swapping two cells; cleared again once the swap is done so
# it leaves no trace in the saved workbook.
$scratchRow = 5000
$scratchCol = 1

$rowPairs = @(
    @(107, 108),
    @(128, 129),
    @(148, 149),
    @(153, 154),
    @(211, 212)
)

$firstCol = 2   # column B
$lastCol  = 30  # column AD

foreach ($pair in $rowPairs) {
    $rowA = $pair[0]
    $rowB = $pair[1]

    foreach ($col in $firstCol..$lastCol) {
        $ws.Cells.Item($scratchRow, $scratchCol).Value2 = $ws.Cells.Item($rowA, $col).Value2
        $ws.Cells.Item($rowA, $col).Value2 = $ws.Cells.Item($rowB, $col).Value2
        $ws.Cells.Item($rowB, $col).Value2 = $ws.Cells.Item($scratchRow, $scratchCol).Value2
    }
}

$ws.Cells.Item($scratchRow, $scratchCol).ClearContents()
